$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Jogos Gerados"

# Update B column values
$ws.Range("B2").Value = 64
$ws.Range("B3").Value = 64
$ws.Range("B4").Value = 54

# Clear the columns/rows that are no longer part of the data (C1:E5)
$ws.Range("C1:E5").Clear()

# Remove row 5 entirely (it was fully removed from the sheet)
$ws.Range("A5:E5").Clear()
